$d = $word.ActiveDocument
$replacements = @(
    @("17+48=65", "40-9=31"),
    @("69+0=69", "10+20=30"),
    @("86-25=61", "83-81=2"),
    @("6+47=53", "50-44=6"),
    @("39+54=93", "17-11=6"),
    @("43-10=33", "45+27=72"),
    @("43-42=1", "37+53=90"),
    @("18+45=63", "94-63=31"),
    @("49+23=72", "94-28=66"),
    @("72-16=56", "69-37=32"),
    @("98-29=69", "99-12=87"),
    @("14+64=78", "40-21=19"),
    @("11+82=93", "69-29=40"),
    @("99-19=80", "34-25=9"),
    @("10+29=39", "43-14=29"),
    @("37-13=24", "59-39=20"),
    @("37+23=60", "34-21=13"),
    @("81-13=68", "62+23=85"),
    @("85-25=60", "28+57=85"),
    @("28-20=8", "31-26=5"),
    @("67-43=24", "6+51=57"),
    @("15+27=42", "35-31=4"),
    @("11+79=90", "61-57=4"),
    @("65-56=9", "12-9=3"),
    @("21+28=49", "22+13=35"),
    @("73-33=40", "86-44=42"),
    @("14+45=59", "96-37=59"),
    @("96-63=33", "52-1=51"),
    @("74-39=35", "23+16=39"),
    @("76-5=71", "11+57=68"),
    @("81-28=53", "46+0=46"),
    @("70-23=47", "76-42=34"),
    @("34+7=41", "31+22=53"),
    @("19+22=41", "65+16=81"),
    @("17-6=11", "25+52=77"),
    @("19+68=87", "13+68=81"),
    @("97-14=83", "64-62=2"),
    @("61-10=51", "0+65=65"),
    @("55-0=55", "92-10=82"),
    @("44-5=39", "99-58=41"),
    @("23+14=37", "71+10=81"),
    @("33+12=45", "39-17=22"),
    @("77-37=40", "30+39=69"),
    @("5+59=64", "62+35=97"),
    @("31+23=54", "63-47=16"),
    @("24+62=86", "66-21=45"),
    @("36+6=42", "76-8=68"),
    @("10-0=10", "34-13=21"),
    @("31+68=99", "62+18=80"),
    @("62-22=40", "40+32=72"),
    @("14+19=33", "50+35=85"),
    @("97-7=90", "58-52=6"),
    @("39+49=88", "80-3=77"),
    @("89-84=5", "77-11=66"),
    @("84+7=91", "17+38=55"),
    @("51-48=3", "46-19=27"),
    @("48-26=22", "58+34=92"),
    @("31-29=2", "84-53=31"),
    @("21+37=58", "72-0=72"),
    @("13+25=38", "41-25=16"),
    @("32+9=41", "95-77=18"),
    @("56+24=80", "30+69=99"),
    @("92-84=8", "44+13=57"),
    @("16+71=87", "17+29=46"),
    @("65-18=47", "64+33=97"),
    @("81-75=6", "86-37=49"),
    @("94-7=87", "74-2=72"),
    @("82-70=12", "28+57=85"),
    @("56-11=45", "93-15=78"),
    @("84+15=99", "74-72=2"),
    @("95-18=77", "60+3=63"),
    @("46+35=81", "34-21=13"),
    @("69-61=8", "0+70=70"),
    @("99-2=97", "34+6=40"),
    @("12+29=41", "17-10=7"),
    @("38+20=58", "76-12=64"),
    @("51-17=34", "15+50=65"),
    @("69-41=28", "64-42=22"),
    @("29-21=8", "66+9=75"),
    @("16-7=9", "26+38=64"),
    @("50+37=87", "38-13=25"),
    @("83-44=39", "94-43=51"),
    @("77-31=46", "90-28=62"),
    @("84-7=77", "93-65=28"),
    @("41+19=60", "85-61=24"),
    @("65-24=41", "87-13=74"),
    @("12-6=6", "46+44=90"),
    @("48+13=61", "36+18=54"),
    @("58-21=37", "60+33=93"),
    @("72-9=63", "64+9=73"),
    @("61+1=62", "4+8=12"),
    @("29+4=33", "87-4=83"),
    @("35+39=74", "47+38=85"),
    @("11+55=66", "59-7=52"),
    @("7+54=61", "89-12=77"),
    @("66+7=73", "6+18=24"),
    @("33-20=13", "54-24=30"),
    @("88-55=33", "71-18=53"),
    @("8+5=13", "95-94=1"),
    @("85-9=76", "76-54=22")
)

$totalReplaced = 0
foreach ($pair in $replacements) {
    $old = $pair[0]
    $new = $pair[1]
    $found = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if ($found) {
        $totalReplaced = $totalReplaced + 1
    } else {
        Write-Output "NOT FOUND: $old"
    }
}
Write-Output "Total replaced: $totalReplaced of $($replacements.Length)"
